$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4703348577022552
$ws.Range("B1").Value = 1.426826953887939
$ws.Range("C1").Value = 4.734130859375
$ws.Range("D1").Value = 1.41790771484375
$ws.Range("E1").Value = 0.8192694187164307
